# Edit script: fills in the blank-form fields in ΕΚΘΕΣΗ_ΣΥΛΛΗΨΗΣ.docx
$d = $word.ActiveDocument

# --- Paragraph 1 (the long "Στην ... του Π.Κ." paragraph): split into 3 runs
#     ("    Στην Θέρμη", a literal tab, and the remaining filled-in sentence)
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Στην  σήμερα την  του μήνα  του έτους  ημέρα εβδομάδας*") {
        $targetPara = $p
        break
    }
}

if ($targetPara -eq $null) {
    throw "Could not locate the target paragraph (blank arrest-report sentence)."
}

$r = $d.Range($targetPara.Range.Start, $targetPara.Range.End - 1)
$xmlFragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">    Στην Θέρμη</w:t></w:r><w:r><w:tab/></w:r><w:r><w:t xml:space="preserve"> σήμερα την 3 του μήνα Νοεμβρίου του έτους 2025 ημέρα εβδομάδας Δευτερα και ώρα 123343δψσδψ ενώπιον εμού του Υ/Α Αθαθααθδσψνβσδκνψ του ψσδκψωδφσωνδφω Θεσσαλονίκης, παρισταμένου και της ωσωσδωσδω της ιδίας υπηρεσίας, που προσλήφθηκε ως Β'' Ανακριτικός Υπάλληλος, οδηγήθηκε στο κατάστημα του Α.Τ. Θέρμης, ο κάτωθι ωδφωω  ωδφωδδφω του ωδφωδφωδφ και της ωδωδφωφδ γεν. ωδωδφωφδδφω στη δφωφωφδ κατ. ωδφωφδωφδ,αριθμός τηλεφώνου ωφδφωφδω, ηλεκτρονικό  ταχυδρομείουωδφωδφωδφ, κάτοχος του υπ αριθμόν ωδωδφ που εκδόθηκε την ωδωδφωδφωφδ από δφσφσδσφδσΑ.Φ.Μ : φσδφδσφδσφδσ, Δ.Ο.Υ : φσδφδσφδσφ, από τον ψαδψδσψ υπηρετών στο ψσδκψωδφσωνδφω, που τον συνέλαβε στις ψσδψδσψδσ και ώρα ψσδψδ στη ψδσψδσψσδψ δσ για παράβαση του/των άρθρων 308 "Σωματικεσ" και 361 ερρωηφ του Π.Κ.</w:t></w:r></w:p>'
$r.InsertXML($xmlFragment)

# --- Paragraph with "Η παρούσα έκθεση άρχισε να συντάσσεται..." : fill in the two time blanks
$find1 = 'Η παρούσα έκθεση άρχισε να συντάσσεται την  ’ ώρα και περατώθηκε την  ώρα.'
$replace1 = 'Η παρούσα έκθεση άρχισε να συντάσσεται την 123343δψσδψ ’ ώρα και περατώθηκε την ψδσψδσψσ ώρα.'
$d.Content.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $replace1, 2)

Write-Output "edit complete"
